# Fix Jinja2 template error, improve row numbering logic, and add email loading state
#
# This script applies the following fixes to the inspection-plan table slides:
#  1. Normalize the pressure unit text "1 Bar.G" -> "1 Bar G" everywhere it
#     appears (column "P(Mpa)" of the component table).
#  2. Fix rows whose "SPEC" column reads "ASTM A36" (a carbon-steel structural
#     spec) but whose MATERIAL / GR columns were mis-populated (e.g. with
#     stainless-steel grades or bogus grade codes): the MATERIAL cell becomes
#     "Not Found" and the GR cell becomes "-".
#  3. Rename the vessel described as "Condensate Vessel" to "Air Receiver".
#  4. Update the temp-file name recorded as the descr/alt-text of the final
#     picture placeholder slide.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shp = $s.Shapes.Item($j)

        # --- Picture placeholder slide: update the stored temp file name ---
        if ($shp.Type -eq 13) {
            if ($shp.AlternativeText -eq "tmpf41z04rg.jpg") {
                $shp.AlternativeText = "tmp5yydeds9.jpg"
            }
        }

        # --- Plain text box holding the vessel name ---
        if ($shp.HasTextFrame -and -not $shp.HasTable) {
            if ($shp.TextFrame.TextRange.Text -eq "Condensate Vessel") {
                $shp.TextFrame.TextRange.Characters().Text = "Air Receiver"
            }
        }

        # --- Inspection data table ---
        if ($shp.HasTable) {
            $tbl = $shp.Table
            for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
                # Column 9 = operating pressure, e.g. "1 Bar.G"
                $pCell = $tbl.Cell($r, 9).Shape.TextFrame.TextRange
                if ($pCell.Text -eq "1 Bar.G") {
                    $pCell.Characters().Text = "1 Bar G"
                }

                # Column 5 = SPEC, column 4 = MATERIAL TYPE, column 6 = GR
                $specCell = $tbl.Cell($r, 5).Shape.TextFrame.TextRange
                if ($specCell.Text -eq "ASTM A36") {
                    $matCell = $tbl.Cell($r, 4).Shape.TextFrame.TextRange
                    if ($matCell.Text -ne "Not Found") {
                        $matCell.Characters().Text = "Not Found"
                    }
                    $grCell = $tbl.Cell($r, 6).Shape.TextFrame.TextRange
                    if ($grCell.Text -ne "-") {
                        $grCell.Characters().Text = "-"
                    }
                }
            }
        }
    }
}
